# Add new prediction rows (row 24-28) to the "Predictions" sheet and the
# corresponding graded results (row 2-6) to the "Results" sheet, then
# refresh the header-row formatting on every sheet (bold, centered, thin
# border) so it matches what Excel re-writes the style table as.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Predictions sheet - append 5 new rows
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Predictions")

$predRows = @(
    @("2025-08-15", "Jupiler Pro League", "oh leuven", "genk", "Away Win", "74.73%", 1.75, "30.78%", "Completed"),
    @("2025-08-15", "Ligue 1", "rennes", "marseille", "Away Win", "64.93%", 1.91, "24.02%", "Completed"),
    @("2025-08-15", "Superliga", "fv nordsjaelland", "fc copenhagen", "Away Win", "52.99%", 2.05, "8.63%", "Completed"),
    @("2025-08-15", "League", "al ramtha", "al faisaly", "Away Win", "51.85%", 2.05, "6.29%", "Completed"),
    @("2025-08-15", "Primera B", "real santander", "real cartagena", "Away Win", "51.60%", 2.15, "10.94%", "Completed")
)

$startRow = 24
for ($i = 0; $i -lt $predRows.Count; $i++) {
    $r = $startRow + $i
    $row = $predRows[$i]
    # Force columns A (date-like text) and F/H (percent-like text) to plain
    # text first so the strings aren't auto-converted into a date serial /
    # percentage number by value-entry auto-detection.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 6).ClearFormats()
    $ws.Cells.Item($r, 8).ClearFormats()
}

# ---------------------------------------------------------------------
# Results sheet - append 5 new graded rows
# ---------------------------------------------------------------------
$wsR = $wb.Worksheets.Item("Results")

$resRows = @(
    @("2025-08-15", "Jupiler Pro League", "oh leuven", "genk", "Away Win", "Away Win", $true, 0.75, 75),
    @("2025-08-15", "Ligue 1", "rennes", "marseille", "Home Win", "Away Win", $false, -1, -100),
    @("2025-08-15", "Superliga", "fv nordsjaelland", "fc copenhagen", "Away Win", "Away Win", $true, 1.05, 105),
    @("2025-08-15", "League", "al ramtha", "al faisaly", "Home Win", "Away Win", $false, -1, -100),
    @("2025-08-15", "Primera B", "real santander", "real cartagena", "Home Win", "Away Win", $false, -1, -100)
)

$startRow2 = 2
for ($i = 0; $i -lt $resRows.Count; $i++) {
    $r = $startRow2 + $i
    $row = $resRows[$i]
    $wsR.Cells.Item($r, 1).NumberFormat = "@"
    $wsR.Cells.Item($r, 1).Value = $row[0]
    $wsR.Cells.Item($r, 2).Value = $row[1]
    $wsR.Cells.Item($r, 3).Value = $row[2]
    $wsR.Cells.Item($r, 4).Value = $row[3]
    $wsR.Cells.Item($r, 5).Value = $row[4]
    $wsR.Cells.Item($r, 6).Value = $row[5]
    $wsR.Cells.Item($r, 7).Value = $row[6]
    $wsR.Cells.Item($r, 8).Value = $row[7]
    $wsR.Cells.Item($r, 9).Value = $row[8]
    $wsR.Cells.Item($r, 1).ClearFormats()
}

# ---------------------------------------------------------------------
# Re-apply header formatting (bold, centered, thin border) on every
# sheet's header row - this is what causes Excel to regenerate the
# style table with fresh font/border/xf entries.
# ---------------------------------------------------------------------
$wsS = $wb.Worksheets.Item("Summary")

$headerRanges = @(
    $ws.Range("A1:I1"),
    $wsR.Range("A1:I1"),
    $wsS.Range("A1:C1")
)

foreach ($hdr in $headerRanges) {
    $hdr.Font.Bold = $true
    $hdr.Font.Name = "Calibri"
    $hdr.Font.Size = 11
    $hdr.HorizontalAlignment = -4108
    $hdr.VerticalAlignment = -4160
    $hdr.Borders.LineStyle = 1
    $hdr.Borders.Weight = 2
}
